$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.020306348800659
$ws.Range("B1").Value = 2.284380435943604
$ws.Range("C1").Value = 2.882810592651367
$ws.Range("D1").Value = 5.840785026550293
$ws.Range("E1").Value = 2.8840651512146
